# 07200005 added review hk
# Rename the "slowthai2" tracklist to "genesisowusu1" and replace the track
# data (45 rows of slowthai tracks -> Genesis Owusu tracks, 15 tracks now
# instead of 14). Sheet1 and Sheet3 hold identical copies of the raw data
# pulled in via web query; Sheet2 is a formatted "pretty print" sheet that
# recalculates automatically off Sheet1 via array formulas.

$wb = $excel.ActiveWorkbook

function Set-TrackData($ws) {
    $ws.Cells.Item(2,2).Value = 'On the Move!'
    $ws.Cells.Item(2,3).Value = 'Kofi Ansah, Andrew Klippel, David Haddad'
    $ws.Cells.Item(2,4).Value = 'Genesis Owusu'
    $ws.Cells.Item(2,5).Value = 0.07430555555555556

    $ws.Cells.Item(3,2).Value = 'The Other Black Dog'
    $ws.Cells.Item(3,3).Value = 'Kofi Ansah, Andrew Klippel, Michael Di Francesco'
    $ws.Cells.Item(3,4).Value = 'Genesis Owusu'
    $ws.Cells.Item(3,5).Value = 0.1826388888888889

    $ws.Cells.Item(4,2).Value = 'Centrefold'
    $ws.Cells.Item(4,3).Value = 'Kofi Ansah, Andrew Klippel, Kieran John Callinan, Di Francesco, Julian Sudek'
    $ws.Cells.Item(4,4).Value = 'Genesis Owusu'
    $ws.Cells.Item(4,5).Value = 0.1361111111111111

    $ws.Cells.Item(5,2).Value = 'Waitin'' on Ya'
    $ws.Cells.Item(5,3).Value = 'Kofi Ansah, Andrew Klippel, Kieran John Callinan, Di Francesco, Julian Sudek'
    $ws.Cells.Item(5,4).Value = 'Genesis Owusu'
    $ws.Cells.Item(5,5).Value = 0.2388888888888889

    $ws.Cells.Item(6,2).Value = 'Don''t Need You'
    $ws.Cells.Item(6,3).Value = 'Kofi Ansah, Andrew Klippel, Kieran John Callinan, Di Francesco, Julian Sudek, David Haddad'
    $ws.Cells.Item(6,4).Value = 'Genesis Owusu'
    $ws.Cells.Item(6,5).Value = 0.12847222222222224

    $ws.Cells.Item(7,2).Value = 'Drown'
    $ws.Cells.Item(7,3).Value = 'Kofi Ansah, Andrew Klippel, Kieran John Callinan, Di Francesco, Julian Sudek, David Haddad'
    $ws.Cells.Item(7,4).Value = 'Genesis Owusu feat. Kirin J. Callinan'
    $ws.Cells.Item(7,5).Value = 0.12291666666666667

    $ws.Cells.Item(8,2).Value = 'Gold Chains'
    $ws.Cells.Item(8,3).Value = 'Kofi Ansah, Andrew Klippel, Kieran John Callinan, Di Francesco, Julian Sudek'
    $ws.Cells.Item(8,4).Value = 'Genesis Owusu'
    $ws.Cells.Item(8,5).Value = 0.15138888888888888

    $ws.Cells.Item(9,2).Value = 'Smiling With No Teeth'
    $ws.Cells.Item(9,3).Value = 'Kofi Ansah, Andrew Klippel, Kieran John Callinan, Di Francesco, Julian Sudek'
    $ws.Cells.Item(9,4).Value = 'Genesis Owusu'
    $ws.Cells.Item(9,5).Value = 0.2076388888888889

    $ws.Cells.Item(10,2).Value = 'I Don''t See Colour'
    $ws.Cells.Item(10,3).Value = 'Kofi Ansah, Andrew Klippel, Kieran John Callinan, Di Francesco, Julian Sudek'
    $ws.Cells.Item(10,4).Value = 'Genesis Owusu'
    $ws.Cells.Item(10,5).Value = 0.12222222222222223

    $ws.Cells.Item(11,2).Value = 'Black Dogs!'
    $ws.Cells.Item(11,3).Value = 'Kofi Ansah, Matt Corby'
    $ws.Cells.Item(11,4).Value = 'Genesis Owusu'
    $ws.Cells.Item(11,5).Value = 0.08333333333333333

    $ws.Cells.Item(12,2).Value = 'Whip Cracker'
    $ws.Cells.Item(12,3).Value = 'Kofi Ansah, Andrew Klippel, Kieran John Callinan, Di Francesco, Julian Sudek'
    $ws.Cells.Item(12,4).Value = 'Genesis Owusu'
    $ws.Cells.Item(12,5).Value = 0.1951388888888889

    $ws.Cells.Item(13,2).Value = 'Easy'
    $ws.Cells.Item(13,3).Value = 'Kofi Ansah, Harvey Sutherland'
    $ws.Cells.Item(13,4).Value = 'Genesis Owusu'
    $ws.Cells.Item(13,5).Value = 0.12638888888888888

    $ws.Cells.Item(14,2).Value = 'A Song About Fishing'
    $ws.Cells.Item(14,3).Value = 'Kofi Ansah, Andrew Klippel, Kieran John Callinan, Di Francesco, Julian Sudek'
    $ws.Cells.Item(14,4).Value = 'Genesis Owusu'
    $ws.Cells.Item(14,5).Value = 0.1423611111111111

    $ws.Cells.Item(15,2).Value = 'No Looking Back'
    $ws.Cells.Item(15,3).Value = 'Kofi Ansah, Andrew Klippel, Kieran John Callinan, Di Francesco, Julian Sudek'
    $ws.Cells.Item(15,4).Value = 'Genesis Owusu'
    $ws.Cells.Item(15,5).Value = 0.15416666666666667

    # Row 16 is a brand new track row (the tracklist grew from 14 to 15
    # tracks), so the "No." column A also needs to be populated.
    $ws.Cells.Item(16,1).Value = 15
    $ws.Cells.Item(16,2).Value = 'Bye Bye'
    $ws.Cells.Item(16,3).Value = 'Kofi Ansah, Andrew Klippel, Kieran John Callinan, Di Francesco, Julian Sudek'
    $ws.Cells.Item(16,4).Value = 'Genesis Owusu'
    $ws.Cells.Item(16,5).Value = 0.17291666666666669
}

$ws1 = $wb.Worksheets.Item("Sheet1")
$ws3 = $wb.Worksheets.Item("Sheet3")

Set-TrackData $ws1
Set-TrackData $ws3

# Column B/D got a bit wider/narrower to fit the new (longer) titles and
# (shorter) performer names.
$ws1.Columns.Item(2).ColumnWidth = 20.43
$ws1.Columns.Item(4).ColumnWidth = 32.93
$ws3.Columns.Item(2).ColumnWidth = 20.43
$ws3.Columns.Item(4).ColumnWidth = 32.93

# The named ranges backing the web query tables get renamed from
# "slowthai2" to "genesisowusu1" and widened by one row (15 -> 16) to cover
# the newly added 15th track.
$n1 = $wb.Names.Item("Sheet1!slowthai2")
$n1.Name = "genesisowusu1"
$n1b = $wb.Names.Item("Sheet1!genesisowusu1")
$n1b.RefersTo = '=Sheet1!$A$1:$E$16'

$n2 = $wb.Names.Item("Sheet3!slowthai2")
$n2.Name = "genesisowusu1"
$n2b = $wb.Names.Item("Sheet3!genesisowusu1")
$n2b.RefersTo = '=Sheet3!$A$1:$E$16'

# Sheet2 is a derived "pretty" report sheet whose cells are all array
# formulas referencing Sheet1 (title/composer/performer/time columns); it
# recalculates automatically from the Set-TrackData writes above. Only its
# manual cell selection (extended by one row for the new 15th track) needs
# to be restored explicitly.
$ws2 = $wb.Worksheets.Item("Sheet2")
$ws2.Activate()
$ws2.Range("A3:K19").Select()
